$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits: add new publisher mappings, update SCW row ---

# 1. Insert "Blackwell Publishing Asia" row at row 7 (alphabetically after "Beacon Press")
$ws.Rows(7).Insert()
$ws.Range("A7").Value = "Blackwell Publishing Asia"
$ws.Range("B7").Value = "Blackwell Publishing Asia"

# 2. Insert "Cornell University Press" row at row 11 (after "Center for Peacebuilding...")
$ws.Rows(11).Insert()
$ws.Range("A11").Value = "Cornell University Press"
$ws.Range("B11").Value = "Cornell University Press"

# 3. Insert combined IOM/ARCM source row at row 23 (after the ARCM row)
$ws.Rows(23).Insert()
$ws.Range("A23").Value = "International Organization on Migration (IOM) and the Asia Research Center on Migration (ARCM), Institute of Asian Studies, Chulalongkorn University."
$ws.Range("B23").Value = "International Organization for Migration; Asian Research Center for Migration; Chulalongkorn University"

# 4. Insert "NUS Press" row at row 27 (after "National Political Publishing House")
$ws.Rows(27).Insert()
$ws.Range("A27").Value = "NUS Press"
$ws.Range("B27").Value = "NUS Press Singapore"
# the row above (National Political Publishing House) carries the red "Bad" style on
# column B; inserting a row copies that formatting down, so reset it here.
$ws.Range("A27:B27").Style = "Normal"
$ws.Rows(27).AutoFit()

# 5. SCW row source got resolved -> update from the placeholder "DO NOT ENTER" text
#    to the real source name, and clear the red "Bad" style that flagged it.
$ws.Range("B30").Value = "Save Cambodia's Wildlife"
$ws.Range("B30").Style = "Normal"
$ws.Rows(30).AutoFit()

# --- Cosmetic / view updates ---
$ws.Columns("A").ColumnWidth = 118.60807291666667
$ws.Columns("B").ColumnWidth = 65.16666666666667
$excel.ActiveWindow.Zoom = 70
$ws.Range("D14").Select()
